$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = "The turn after this move is used, the Pokémon''s action is skipped so it can recharge."
$ws.Range("D15").Value = "This move ignores the target''s []{move:substitute}."
